$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 7: "Commit hash" row ---
# The original paragraph has "Commit hash" followed by a collapsed
# _GoBack bookmark. That bookmark needs to move to the new
# "Compare and pull request button" row below, so clear this cell
# completely (which also drops the bookmark) and retype the text.
$cell7 = $t.Cell(7, 1)
$cell7.Range.Delete()
$t.Cell(7, 1).Range.Text = "Commit hash"

# --- Row 8: "Branch" / description ---
$cell8a = $t.Cell(8, 1)
$cell8a.Range.Text = "Branch"
$t.Cell(8, 1).Range.Font.Bold = 1

$cell8b = $t.Cell(8, 2)
$cell8b.Range.Text = "Allows " + [char]0x201C + "experiments" + [char]0x201D + " to the main project.  If the experiment is successful, can merge back to the master."

# --- Row 9: "Graphs -> Network" / description ---
$cell9a = $t.Cell(9, 1)
$cell9a.Range.Text = "Graphs -> Network"
$cell9aFresh = $t.Cell(9, 1)
$start9 = $cell9aFresh.Range.Start
$d.Range($start9, $start9 + 6).Font.Bold = 1
$d.Range($start9 + 10, $start9 + 17).Font.Bold = 1

$cell9b = $t.Cell(9, 2)
$cell9b.Range.Text = "View the state of the repository"

# --- Row 10: "Pull request" / description ---
$cell10a = $t.Cell(10, 1)
$cell10a.Range.Text = "Pull request"
$cell10aFresh = $t.Cell(10, 1)
$start10 = $cell10aFresh.Range.Start
$d.Range($start10, $start10 + 4).Font.Bold = 1

$cell10b = $t.Cell(10, 2)
$cell10b.Range.Text = "Import changes from one branch to another branch."

# --- Row 11: "Merge request" / (empty) ---
$cell11a = $t.Cell(11, 1)
$cell11a.Range.Text = "Merge request"
$cell11aFresh = $t.Cell(11, 1)
$start11 = $cell11aFresh.Range.Start
$d.Range($start11, $start11 + 5).Font.Bold = 1

# --- Row 12: "Compare and pull request button" / (unchanged) ---
$cell12a = $t.Cell(12, 1)
$cell12a.Range.Text = "Compare and pull request button"
$cell12aFresh = $t.Cell(12, 1)
$start12 = $cell12aFresh.Range.Start
$bmRange = $d.Range($start12 + 25, $start12 + 25)
$d.Bookmarks.Add("_GoBack", $bmRange)
